# Expand the changelog bullet:
#   "Hiển thị lượng quân khi xem nhà trại lính"
# into:
#   "Hiển thị Content: lượng lính hiện tại, lính mở khóa, công trình mở khóa,..
#    khi nhấn xem thông tin hoặc nâng cấp công trình"
# while keeping the pre-existing _GoBack bookmark sitting right before
# " nâng cấp công trình" (its original position relative to the sentence
# end), split across the same run boundaries the source diff shows.

$d = $word.ActiveDocument

$oldText  = "Hiển thị lượng quân khi xem nhà trại lính"

# The four textual segments the final paragraph is made of (in order).
# part3/part4 are split by the relocated _GoBack bookmark.
$part1 = "Hiển thị "
$part2 = "Content: lượ"
$part3 = "ng lính hiện tại, lính mở khóa, công trình mở khóa,.. khi nhấn xem thông tin hoặc"
$part4 = " nâng cấp công trình"
$newText = $part1 + $part2 + $part3 + $part4

# Drop the existing bookmark first so the text substitution below isn't
# constrained by its current anchor; it gets re-created afterwards at the
# right spot.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Replace the whole sentence in one shot (this lands as a single run).
$r = $d.Content
$r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
$runStart = $r.Start

# Offsets (relative to the document) of the three split points inside the
# replaced run.
$split1 = $runStart + $part1.Length                 # between part1 / part2
$split2 = $split1   + $part2.Length                 # between part2 / part3
$split3 = $split2   + $part3.Length                 # between part3 / part4 (bookmark goes here)

# Dropping a zero-length bookmark at a position forces Word to cleanly
# split the run there without leaving stray run formatting behind, so use
# two throw-away bookmarks to get the part1|part2|part3 boundaries, then
# recreate _GoBack at the part3|part4 boundary (which both splits the run
# there and restores the bookmark in its rightful place).
$d.Bookmarks.Add("TempSplitA", $d.Range($split1, $split1))
$d.Bookmarks.Add("TempSplitB", $d.Range($split2, $split2))
$d.Bookmarks.Add("_GoBack", $d.Range($split3, $split3))

$d.Bookmarks("TempSplitA").Delete()
$d.Bookmarks("TempSplitB").Delete()
